$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new column header
$ws.Range("I1").Value = "Lifeyearlost70_CI"

# New "Lifeyearlost70_CI" values for each data row (row 14 / Sweden 2021 stays blank,
# matching the rest of that row which has no data either)
$ws.Range("I2").Value = "-38 to 136"
$ws.Range("I3").Value = "109 to 361"
$ws.Range("I4").Value = "834 to 965"
$ws.Range("I5").Value = "945 to 1189"
$ws.Range("I6").Value = "9132 to 10926"
$ws.Range("I7").Value = "-51 to 6"
$ws.Range("I8").Value = "-13 to 30"
$ws.Range("I9").Value = "-383 to 172"
$ws.Range("I10").Value = "-3 to 15"
$ws.Range("I11").Value = "-6 to 13"
$ws.Range("I12").Value = "39 to 145"
$ws.Range("I13").Value = "-12 to 1"
$ws.Range("I14").Value = ""
$ws.Range("I15").Value = "-32 to 41"
